$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$headers = @("Name","Team","ERA","G","W","L","SV","HLD","WPCT","IP","H","HR","BB","HBP","SO","R","ER","WHIP")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data row (row 2) : 폰세 / 한화 pitcher stats ----
$values = @("폰세","한화",1.8,13,9,0,0,0,1,85,52,4,20,1,112,18,17,0.85)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}

# ---- Formatting ----
# Name / Team / ERA columns (A:C) keep the base font.
# Stat columns (D:R) use the 돋움 font like the rest of the rebuilt table.
$statHeader = $ws.Range("D1:R1")
$statHeader.Font.Name = "돋움"
$statHeader.Font.Size = 12
$statHeader.Font.ThemeColor = 1

$statData = $ws.Range("D2:R2")
$statData.Font.Name = "돋움"
$statData.Font.Size = 12
$statData.Font.ThemeColor = 1

$ws.Rows.Item(1).Font.Size = 12
$ws.Rows.Item(2).Font.Size = 12

# ---- Sheet view bookkeeping to mirror the saved workbook state ----
$ws.Range("I6").Select()

Write-Host "applied pitcher stats table"
